$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'25.788.08"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.53%  '

$ws.Range('D3').Value = "'1.637.37"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.25%  '

$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = "'215.50"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '

$ws.Range('D6').Value = "'0.5058"
$ws.Range('D6').Style = 'Normal'

$ws.Range('D7').Value = "'1.002"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').Value = "'0.2581"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.22%  '

$ws.Range('D9').Value = "'0.06418"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.22%  '

$ws.Range('E10').Value = '  +4.38%  '

$ws.Range('D11').Value = "'0.07787"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.19%  '

$ws.Range('E12').Value = '  -0.28%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = "'1.864.39"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.33%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = "'1.636.38"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.02%  '

$ws.Range('D15').Value = "'0.5612"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.81%  '

$ws.Range('D16').Value = "'0.0₅7637"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.10%  '

$ws.Range('D17').Value = "'63.27"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.06%  '

$ws.Range('D18').Value = "'25.812.36"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.53%  '

$ws.Range('E19').Value = '  -0.05%  '

$ws.Range('D20').Value = "'192.96"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.07%  '

$ws.Range('D21').Value = "'4.379"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.93%  '

$ws.Range('D22').Value = "'9.913"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.42%  '

$ws.Range('D23').Value = "'6.128"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.28%  '

$ws.Range('D24').Value = "'1.002"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.06%  '

$ws.Range('D25').Value = "'1.798"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.82%  '

$ws.Range('D26').Value = "'140.79"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.80%  '

$ws.Range('E27').Value = '  -1.79%  '

$ws.Range('D28').Value = "'6.811"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.89%  '

$ws.Range('D29').Value = "'15.56"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.08%  '

$ws.Range('E30').Value = '  +0.31%  '

$ws.Range('D31').Value = "'0.04952"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.50%  '

$ws.Range('D32').Value = "'3.282"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.18%  '

$ws.Range('D33').Value = "'3.232"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.28%  '

$ws.Range('D34').Value = "'1.569"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.78%  '

$ws.Range('D35').Value = "'2.383"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.45%  '

$ws.Range('D36').Value = "'0.9033"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.68%  '

$ws.Range('D37').Value = "'0.5579"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.30%  '

$ws.Range('D38').Value = "'2.573"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.37%  '

$ws.Range('D39').Value = "'1.131.41"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.25%  '

$ws.Range('D40').Value = "'0.01568"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.75%  '

$ws.Range('D41').Value = "'0.9952"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.62%  '

$ws.Range('D42').Value = "'5.486"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.90%  '

$ws.Range('D43').Value = "'0.8024"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.83%  '

$ws.Range('D44').Value = "'98.82"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.28%  '

$ws.Range('D45').Value = "'1.775.76"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.39%  '

$ws.Range('E46').Value = '  -5.50%  '

$ws.Range('D47').Value = "'55.59"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.55%  '

$ws.Range('D48').Value = "'0.4270"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.95%  '

$ws.Range('D49').Value = "'7.759"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.42%  '

$ws.Range('D50').Value = "'0.05032"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.92%  '

$ws.Range('D51').Value = "'0.9987"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.44%  '
